$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 4.915327667582776
$ws.Cells.Item(2, 4).Value = 8.60625054328043
$ws.Cells.Item(2, 5).Value = 13.06961900017016
$ws.Cells.Item(2, 6).Value = 34.2322369635045
$ws.Cells.Item(2, 7).Value = 3.666014162855763
$ws.Cells.Item(2, 10).Value = 9.929445973404787
$ws.Cells.Item(2, 11).Value = 15.35771083111988
$ws.Cells.Item(2, 13).Value = 17.63982092207302
$ws.Cells.Item(2, 15).Value = 26.06952759829517
$ws.Cells.Item(3, 3).Value = 4.751136923475149
$ws.Cells.Item(3, 4).Value = 8.583700291697944
$ws.Cells.Item(3, 5).Value = 13.09279494133811
$ws.Cells.Item(3, 6).Value = 34.36639315128756
$ws.Cells.Item(3, 7).Value = 3.66843235246811
$ws.Cells.Item(3, 10).Value = 9.961639534516067
$ws.Cells.Item(3, 11).Value = 14.84199935339796
$ws.Cells.Item(3, 13).Value = 17.42979383864191
$ws.Cells.Item(3, 15).Value = 26.20177038389949
$ws.Cells.Item(4, 3).Value = 4.648899802561806
$ws.Cells.Item(4, 4).Value = 8.570756156284787
$ws.Cells.Item(4, 5).Value = 13.10913539096495
$ws.Cells.Item(4, 6).Value = 34.45889491948046
$ws.Cells.Item(4, 7).Value = 3.669994763443318
$ws.Cells.Item(4, 10).Value = 9.982631862672298
$ws.Cells.Item(4, 11).Value = 14.51694248570413
$ws.Cells.Item(4, 13).Value = 17.3016017346342
$ws.Cells.Item(4, 15).Value = 26.29027466843648
$ws.Cells.Item(5, 3).Value = 4.606954362739823
$ws.Cells.Item(5, 4).Value = 8.565711138501667
$ws.Cells.Item(5, 5).Value = 13.11632476446715
$ws.Cells.Item(5, 6).Value = 34.49912585639795
$ws.Cells.Item(5, 7).Value = 3.670651044414726
$ws.Cells.Item(5, 10).Value = 9.991495006474928
$ws.Cells.Item(5, 11).Value = 14.38255056096052
$ws.Cells.Item(5, 13).Value = 17.24960546646446
$ws.Cells.Item(5, 15).Value = 26.3281712595824
$ws.Cells.Item(6, 3).Value = 4.599974560586289
$ws.Cells.Item(6, 4).Value = 8.564887383960814
$ws.Cells.Item(6, 5).Value = 13.11755059007074
$ws.Cells.Item(6, 6).Value = 34.50595898372826
$ws.Cells.Item(6, 7).Value = 3.67076120419943
$ws.Cells.Item(6, 10).Value = 9.992985377488059
$ws.Cells.Item(6, 11).Value = 14.3601241579516
$ws.Cells.Item(6, 13).Value = 17.24098767782301
$ws.Cells.Item(6, 15).Value = 26.33457430655873
$ws.Cells.Item(7, 3).Value = 4.648335154385606
$ws.Cells.Item(7, 4).Value = 8.570687182943185
$ws.Cells.Item(7, 5).Value = 13.10923020161433
$ws.Cells.Item(7, 6).Value = 34.4594272363326
$ws.Cells.Item(7, 7).Value = 3.67000353489026
$ws.Cells.Item(7, 10).Value = 9.982750143990698
$ws.Cells.Item(7, 11).Value = 14.51513758036532
$ws.Cells.Item(7, 13).Value = 17.30089944406624
$ws.Cells.Item(7, 15).Value = 26.29077835323674
$ws.Cells.Item(8, 3).Value = 4.859053332766459
$ws.Cells.Item(8, 4).Value = 8.59828992348654
$ws.Cells.Item(8, 5).Value = 13.07717199341025
$ws.Cells.Item(8, 6).Value = 34.27638494106843
$ws.Cells.Item(8, 7).Value = 3.666831880811502
$ws.Cells.Item(8, 10).Value = 9.940292277235246
$ws.Cells.Item(8, 11).Value = 15.18174471946128
$ws.Cells.Item(8, 13).Value = 17.56727515602407
$ws.Cells.Item(8, 15).Value = 26.11360424016588
$ws.Cells.Item(9, 3).Value = 5.257974874182158
$ws.Cells.Item(9, 4).Value = 8.659434449712853
$ws.Cells.Item(9, 5).Value = 13.03105555855071
$ws.Cells.Item(9, 6).Value = 33.99830394403976
$ws.Cells.Item(9, 7).Value = 3.66122531375461
$ws.Cells.Item(9, 10).Value = 9.866733574309698
$ws.Cells.Item(9, 11).Value = 16.41494500237435
$ws.Cells.Item(9, 13).Value = 18.09348284959993
$ws.Cells.Item(9, 15).Value = 25.82445113628345
$ws.Cells.Item(10, 3).Value = 5.538741265894954
$ws.Cells.Item(10, 4).Value = 8.70843894012935
$ws.Cells.Item(10, 5).Value = 13.0073892574817
$ws.Cells.Item(10, 6).Value = 33.84394338577175
$ws.Cells.Item(10, 7).Value = 3.657475764567983
$ws.Cells.Item(10, 10).Value = 9.818573994055239
$ws.Cells.Item(10, 11).Value = 17.26723596142308
$ws.Cells.Item(10, 13).Value = 18.47946372814801
$ws.Cells.Item(10, 15).Value = 25.64796422343681
$ws.Cells.Item(11, 3).Value = 5.6631488449933
$ws.Cells.Item(11, 4).Value = 8.731574031753036
$ws.Cells.Item(11, 5).Value = 12.9988405306453
$ws.Cells.Item(11, 6).Value = 33.78469535135773
$ws.Cells.Item(11, 7).Value = 3.655849372084386
$ws.Cells.Item(11, 10).Value = 9.797936363405039
$ws.Cells.Item(11, 11).Value = 17.6417855248859
$ws.Cells.Item(11, 13).Value = 18.65430161004878
$ws.Cells.Item(11, 15).Value = 25.57557279927125
$ws.Cells.Item(12, 3).Value = 5.709736007642201
$ws.Cells.Item(12, 4).Value = 8.740451838274899
$ws.Cells.Item(12, 5).Value = 12.99592197984225
$ws.Cells.Item(12, 6).Value = 33.76384715286247
$ws.Cells.Item(12, 7).Value = 3.655244835473668
$ws.Cells.Item(12, 10).Value = 9.790303666889072
$ws.Cells.Item(12, 11).Value = 17.78161826660041
$ws.Cells.Item(12, 13).Value = 18.72035013795032
$ws.Cells.Item(12, 15).Value = 25.54930254362755
$ws.Cells.Item(13, 3).Value = 5.699726589240771
$ws.Cells.Item(13, 4).Value = 8.738534705134031
$ws.Cells.Item(13, 5).Value = 12.99653637272266
$ws.Cells.Item(13, 6).Value = 33.76826641570745
$ws.Cells.Item(13, 7).Value = 3.655374529717502
$ws.Cells.Item(13, 10).Value = 9.791939401202695
$ws.Cells.Item(13, 11).Value = 17.75159336829838
$ws.Cells.Item(13, 13).Value = 18.70613327239433
$ws.Cells.Item(13, 15).Value = 25.55490936384382
$ws.Cells.Item(14, 3).Value = 5.666992334975352
$ws.Cells.Item(14, 4).Value = 8.732302092511583
$ws.Cells.Item(14, 5).Value = 12.99859403440574
$ws.Cells.Item(14, 6).Value = 33.78294828106539
$ws.Cells.Item(14, 7).Value = 3.655799409505498
$ws.Cells.Item(14, 10).Value = 9.797304764782902
$ws.Cells.Item(14, 11).Value = 17.65333032713171
$ws.Cells.Item(14, 13).Value = 18.65973893431823
$ws.Cells.Item(14, 15).Value = 25.57338858182312
$ws.Cells.Item(15, 3).Value = 5.646872216529079
$ws.Cells.Item(15, 4).Value = 8.728499555810174
$ws.Cells.Item(15, 5).Value = 12.99989590506498
$ws.Cells.Item(15, 6).Value = 33.79214840815364
$ws.Cells.Item(15, 7).Value = 3.656061136060582
$ws.Cells.Item(15, 10).Value = 9.800614938764408
$ws.Cells.Item(15, 11).Value = 17.59287778823541
$ws.Cells.Item(15, 13).Value = 18.63129888051856
$ws.Cells.Item(15, 15).Value = 25.58485669244501
$ws.Cells.Item(16, 3).Value = 5.530540131447559
$ws.Cells.Item(16, 4).Value = 8.706943614760608
$ws.Cells.Item(16, 5).Value = 13.00799254050631
$ws.Cells.Item(16, 6).Value = 33.84803703060333
$ws.Cells.Item(16, 7).Value = 3.657583643665869
$ws.Cells.Item(16, 10).Value = 9.819948244847938
$ws.Cells.Item(16, 11).Value = 17.24248384870513
$ws.Cells.Item(16, 13).Value = 18.46801821251968
$ws.Cells.Item(16, 15).Value = 25.65285467161075
$ws.Cells.Item(17, 3).Value = 5.458290344937025
$ws.Cells.Item(17, 4).Value = 8.693932683013465
$ws.Cells.Item(17, 5).Value = 13.01352735169284
$ws.Cells.Item(17, 6).Value = 33.88514027062617
$ws.Cells.Item(17, 7).Value = 3.658537919777912
$ws.Cells.Item(17, 10).Value = 9.83213372621802
$ws.Cells.Item(17, 11).Value = 17.02407609160443
$ws.Cells.Item(17, 13).Value = 18.36762259143918
$ws.Cells.Item(17, 15).Value = 25.69659643978018
$ws.Cells.Item(18, 3).Value = 5.416424646589278
$ws.Cells.Item(18, 4).Value = 8.686528714195168
$ws.Cells.Item(18, 5).Value = 13.01691953652171
$ws.Cells.Item(18, 6).Value = 33.90751318145999
$ws.Cells.Item(18, 7).Value = 3.659094261443188
$ws.Cells.Item(18, 10).Value = 9.83926208988758
$ws.Cells.Item(18, 11).Value = 16.89722046982205
$ws.Cells.Item(18, 13).Value = 18.30981040085582
$ws.Cells.Item(18, 15).Value = 25.7224980476357
$ws.Cells.Item(19, 3).Value = 5.40219798760027
$ws.Cells.Item(19, 4).Value = 8.684035648327571
$ws.Cells.Item(19, 5).Value = 13.01810392133907
$ws.Cells.Item(19, 6).Value = 33.9152652138789
$ws.Cells.Item(19, 7).Value = 3.659283913513665
$ws.Cells.Item(19, 10).Value = 9.841696183284302
$ws.Cells.Item(19, 11).Value = 16.85406104325749
$ws.Cells.Item(19, 13).Value = 18.29022621656737
$ws.Cells.Item(19, 15).Value = 25.73139517504723
$ws.Cells.Item(20, 3).Value = 5.466013848398143
$ws.Cells.Item(20, 4).Value = 8.695309512707636
$ws.Cells.Item(20, 5).Value = 13.01291656202316
$ws.Cells.Item(20, 6).Value = 33.88108367313269
$ws.Cells.Item(20, 7).Value = 3.658435563054617
$ws.Cells.Item(20, 10).Value = 9.830824185681402
$ws.Cells.Item(20, 11).Value = 17.04745443455666
$ws.Cells.Item(20, 13).Value = 18.378317216698
$ws.Cells.Item(20, 15).Value = 25.69186315157459
$ws.Cells.Item(21, 3).Value = 5.676621726660722
$ws.Cells.Item(21, 4).Value = 8.734129618790117
$ws.Cells.Item(21, 5).Value = 12.99798100250421
$ws.Cells.Item(21, 6).Value = 33.77859269290246
$ws.Cells.Item(21, 7).Value = 3.655674304642679
$ws.Cells.Item(21, 10).Value = 9.795723881946163
$ws.Cells.Item(21, 11).Value = 17.68224767121075
$ws.Cells.Item(21, 13).Value = 18.6733707886591
$ws.Cells.Item(21, 15).Value = 25.56792971066876
$ws.Cells.Item(22, 3).Value = 5.811198025585145
$ws.Cells.Item(22, 4).Value = 8.760181050532125
$ws.Cells.Item(22, 5).Value = 12.99007704583275
$ws.Cells.Item(22, 6).Value = 33.72086777027088
$ws.Cells.Item(22, 7).Value = 3.653935750442916
$ws.Cells.Item(22, 10).Value = 9.773846383536757
$ws.Cells.Item(22, 11).Value = 18.08541682052683
$ws.Cells.Item(22, 13).Value = 18.86525694398528
$ws.Cells.Item(22, 15).Value = 25.49359673977315
$ws.Cells.Item(23, 3).Value = 5.739666929392205
$ws.Cells.Item(23, 4).Value = 8.74621606584245
$ws.Cells.Item(23, 5).Value = 12.99412567009149
$ws.Cells.Item(23, 6).Value = 33.75082628426727
$ws.Cells.Item(23, 7).Value = 3.654857622148429
$ws.Cells.Item(23, 10).Value = 9.785425705200046
$ws.Cells.Item(23, 11).Value = 17.87134134095521
$ws.Cells.Item(23, 13).Value = 18.76294711259991
$ws.Cells.Item(23, 15).Value = 25.53265728040844
$ws.Cells.Item(24, 3).Value = 5.462523071606399
$ws.Cells.Item(24, 4).Value = 8.694686810388639
$ws.Cells.Item(24, 5).Value = 13.01319204543159
$ws.Cells.Item(24, 6).Value = 33.88291441652004
$ws.Cells.Item(24, 7).Value = 3.658481814491047
$ws.Cells.Item(24, 10).Value = 9.831415846507435
$ws.Cells.Item(24, 11).Value = 17.03688909570501
$ws.Cells.Item(24, 13).Value = 18.37348246363489
$ws.Cells.Item(24, 15).Value = 25.69400072338112
$ws.Cells.Item(25, 3).Value = 5.151981794312768
$ws.Cells.Item(25, 4).Value = 8.642161711870679
$ws.Cells.Item(25, 5).Value = 13.04173735028324
$ws.Cells.Item(25, 6).Value = 34.06480637533497
$ws.Cells.Item(25, 7).Value = 3.662676837090531
$ws.Cells.Item(25, 10).Value = 9.885597728209659
$ws.Cells.Item(25, 11).Value = 16.09020723062935
$ws.Cells.Item(25, 13).Value = 17.95102616608301
$ws.Cells.Item(25, 15).Value = 25.89639126278806
